$d = $word.ActiveDocument

# --- Step 1: remove the existing "_GoBack" bookmark (currently sitting at
#     the end of the YouTube-link paragraph, right after the hyperlink run).
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBookmark = $d.Bookmarks("_GoBack")
    $oldBookmark.Delete()
}

# --- Step 2: insert a new, otherwise-empty paragraph directly after the
#     "Link to video" paragraph. Assigning the paragraph-mark character to
#     the collapsed end-of-paragraph range duplicates that paragraph's
#     pPr/rPr (spacing, indent, bold Times New Roman) without adding any
#     stray run.
$findRange = $d.Content
$findRange.Find.Execute("Link to video", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$linkPara = $findRange.Paragraphs(1)
$insertionPoint = $linkPara.Range
$insertionPoint.Collapse(0)
$insertionPoint.Text = "`r"

# --- Step 3: re-create the "_GoBack" bookmark, collapsed, inside the new
#     empty paragraph we just inserted.
$findRange2 = $d.Content
$findRange2.Find.Execute("Link to video", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$linkPara2 = $findRange2.Paragraphs(1)
$newPara = $linkPara2.Next()
$newParaRange = $newPara.Range
$newParaRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $newParaRange)
